$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-31 Thursday" "2023-09-01 Friday"

Replace-Text "27×82=" "88×81="
Replace-Text "74×36=" "94×80="
Replace-Text "69×61=" "43×25="
Replace-Text "25×70=" "80×14="
Replace-Text "26×60=" "59×51="

Replace-Text "71×60=" "39×60="
Replace-Text "33×36=" "61×31="
Replace-Text "97×44=" "53×55="
Replace-Text "31×23=" "76×57="
Replace-Text "61×55=" "66×58="

Replace-Text "44×18=" "18×35="
Replace-Text "68×63=" "67×77="
Replace-Text "80×46=" "67×36="
Replace-Text "98×39=" "53×64="
Replace-Text "36×30=" "14×15="

Replace-Text "99×75=" "54×17="
Replace-Text "12×37=" "42×43="
Replace-Text "72×49=" "91×87="
Replace-Text "43×78=" "44×90="
Replace-Text "81×85=" "69×81="

Replace-Text "59×45=" "12×34="
Replace-Text "99×63=" "24×61="
Replace-Text "45×97=" "80×39="
Replace-Text "58×27=" "39×93="
Replace-Text "36×77=" "61×60="
